$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13 for the "Docentes responsaveis" name, shifting old rows 13-21 down to 14-22
$ws.Rows.Item(13).Insert()

# Copy the B/C cell formatting (style) from row 14 (shifted former row 13) onto new row 13
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 13 column A should stay empty (no cell at all)
$ws.Range("A13").Clear()

$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"
$ws.Range("B2").Value = "8800006"
$ws.Range("C2").Value = "8800006"
$ws.Range("A3").Value = "Nome:"
$ws.Range("B3").Value = " Projeto de Engenharia I"
$ws.Range("C3").Value = " Projeto de Engenharia I"
$ws.Range("A4").Value = "Name:"
$ws.Range("B4").Value = "Engineering Project I"
$ws.Range("C4").Value = "Engineering Project I"
$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("B5").Value = "2"
$ws.Range("C5").Value = "2"
$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("B6").Value = "2"
$ws.Range("C6").Value = "2"
$ws.Range("A7").Value = "Carga horária:"
$ws.Range("B7").Value = "90 h"
$ws.Range("C7").Value = "90 h"
$ws.Range("A8").Value = "Ativação:"
$ws.Range("B8").Value = "01/01/2015"
$ws.Range("C8").Value = "01/01/2015"
$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("B9").Value = "EP-9"
$ws.Range("C9").Value = "EP-9"
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "Desenvolver um projeto sobre tema de Engenharia de Produção, similar a situações que os alunos irão encontrar na vida real no efetivo exercício de sua profissão, `nAplicar e integrar conhecimentos adquiridos em demais disciplinas de seu curso`nDesenvolver competências técnicas, as relacionadas ao projeto em si, bem como competências transversais (habilidades e atitudes), num ambiente de aprendizagem baseado em PBL (Project-Baed Learning)."
$ws.Range("C10").Value = "Desenvolver um projeto sobre tema de Engenharia de Produção, similar a situações que os alunos irão encontrar na vida real no efetivo exercício de sua profissão, `nAplicar e integrar conhecimentos adquiridos em demais disciplinas de seu curso`nDesenvolver competências técnicas, as relacionadas ao projeto em si, bem como competências transversais (habilidades e atitudes), num ambiente de aprendizagem baseado em PBL (Project-Baed Learning)."
$ws.Range("A11").Value = "Objectives:"
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("B13").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C13").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "Tópicos que abordem o tema do projeto de seu planejamento a execução."
$ws.Range("C14").Value = "Tópicos que abordem o tema do projeto de seu planejamento a execução."
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "Noções de Gestão de Projetos`nOrganização do tempo: dimensão pessoal;`nTécnicas para a realização de apresentações;`nNoções de Aprendizagem Baseada em Projetos`nTrabalho em Grupo, Equipes e times. `nPostura e Ética Profissional`nTécnicas para redação de relatório técnico;`nTutoria de projetos.`nAssuntos Técnicos específicos relacionados com o tema do projeto."
$ws.Range("C16").Value = "Noções de Gestão de Projetos`nOrganização do tempo: dimensão pessoal;`nTécnicas para a realização de apresentações;`nNoções de Aprendizagem Baseada em Projetos`nTrabalho em Grupo, Equipes e times. `nPostura e Ética Profissional`nTécnicas para redação de relatório técnico;`nTutoria de projetos.`nAssuntos Técnicos específicos relacionados com o tema do projeto."
$ws.Range("A17").Value = "Syllabus:"
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras.`n`nOs alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Produção, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão. `nCada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.`nAs aulas ocorrerão: 1) através de uma reunião da equipe de trabalho para tratar do projeto, e  2) palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores  ou profissionais de empresas."
$ws.Range("C19").Value = "O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras.`n`nOs alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Produção, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão. `nCada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.`nAs aulas ocorrerão: 1) através de uma reunião da equipe de trabalho para tratar do projeto, e  2) palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores  ou profissionais de empresas."
$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "A nota será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.`nO detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na coordenação da disciplina."
$ws.Range("C20").Value = "A nota será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.`nO detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na coordenação da disciplina."
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "Não há recuperação"
$ws.Range("C21").Value = "Não há recuperação"
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "Artigos sobre metodologias ativas de aprendizagem e  Project Based Learning.`nLivros e Artigos científicos relacionados com o tema do projeto."
$ws.Range("C22").Value = "Artigos sobre metodologias ativas de aprendizagem e  Project Based Learning.`nLivros e Artigos científicos relacionados com o tema do projeto."
